$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '96.588.03'
Set-TextValue $ws.Range('E2') '  -2.13%  '

Set-TextValue $ws.Range('D3') '3.305.93'
Set-TextValue $ws.Range('E3') '  -4.85%  '

Set-TextValue $ws.Range('E4') '  +0.11%  '

Set-TextValue $ws.Range('D5') '245.95'
Set-TextValue $ws.Range('E5') '  -6.31%  '

Set-TextValue $ws.Range('D6') '648.09'
Set-TextValue $ws.Range('E6') '  -4.11%  '

Set-TextValue $ws.Range('D7') '1.33'
Set-TextValue $ws.Range('E7') '  -15.46%  '

Set-TextValue $ws.Range('D8') '0.408'
Set-TextValue $ws.Range('E8') '  -11.55%  '

Set-TextValue $ws.Range('E9') '  +0.13%  '

Set-TextValue $ws.Range('D10') '0.963'
Set-TextValue $ws.Range('E10') '  -14.91%  '

Set-TextValue $ws.Range('D11') '3.302.37'
Set-TextValue $ws.Range('E11') '  -4.90%  '

Set-TextValue $ws.Range('E12') '  -7.35%  '

Set-TextValue $ws.Range('D13') '39.44'
Set-TextValue $ws.Range('E13') '  -9.34%  '

Set-TextValue $ws.Range('D14') '96.528.14'
Set-TextValue $ws.Range('E14') '  -1.81%  '

Set-TextValue $ws.Range('D15') '5.92'
Set-TextValue $ws.Range('E15') '  -6.13%  '

Set-TextValue $ws.Range('B16') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D16') '3.927.76'
Set-TextValue $ws.Range('E16') '  -4.56%  '

Set-TextValue $ws.Range('B17') 'ShibaInu'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D17') '0.0000247'
Set-TextValue $ws.Range('E17') '  -9.87%  '

Set-TextValue $ws.Range('D18') '8.44'
Set-TextValue $ws.Range('E18') '  -4.38%  '

Set-TextValue $ws.Range('D19') '3.313.07'
Set-TextValue $ws.Range('E19') '  -4.81%  '

Set-TextValue $ws.Range('D20') '16.60'
Set-TextValue $ws.Range('E20') '  -6.87%  '

Set-TextValue $ws.Range('B21') 'BitcoinCash'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D21') '492.86'
Set-TextValue $ws.Range('E21') '  -7.59%  '

Set-TextValue $ws.Range('B22') 'SuiNetwork'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D22') '3.31'
Set-TextValue $ws.Range('E22') '  -8.82%  '

Set-TextValue $ws.Range('D23') '10.26'
Set-TextValue $ws.Range('E23') '  -7.35%  '

Set-TextValue $ws.Range('D24') '0.450'
Set-TextValue $ws.Range('E24') '  -8.48%  '

Set-TextValue $ws.Range('D25') '0.0000195'
Set-TextValue $ws.Range('E25') '  -10.47%  '

Set-TextValue $ws.Range('D26') '6.42'
Set-TextValue $ws.Range('E26') '  -1.41%  '

Set-TextValue $ws.Range('D27') '93.53'
Set-TextValue $ws.Range('E27') '  -10.05%  '

Set-TextValue $ws.Range('D28') '11.86'
Set-TextValue $ws.Range('E28') '  -10.05%  '

Set-TextValue $ws.Range('D29') '3.492.02'
Set-TextValue $ws.Range('E29') '  -4.55%  '

Set-TextValue $ws.Range('E30') '  +0.08%  '

Set-TextValue $ws.Range('E31') '  -9.40%  '

Set-TextValue $ws.Range('D32') '10.61'
Set-TextValue $ws.Range('E32') '  -10.35%  '

Set-TextValue $ws.Range('D33') '0.183'
Set-TextValue $ws.Range('E33') '  -7.51%  '

Set-TextValue $ws.Range('E34') '  +7.22%  '

Set-TextValue $ws.Range('D35') '0.995'
Set-TextValue $ws.Range('E35') '  -0.39%  '

Set-TextValue $ws.Range('D36') '0.534'
Set-TextValue $ws.Range('E36') '  -11.19%  '

Set-TextValue $ws.Range('D37') '27.65'
Set-TextValue $ws.Range('E37') '  -9.81%  '

Set-TextValue $ws.Range('D38') '1.45'
Set-TextValue $ws.Range('E38') '  +0.50%  '

Set-TextValue $ws.Range('D39') '7.43'
Set-TextValue $ws.Range('E39') '  -9.16%  '

Set-TextValue $ws.Range('D41') '0.148'
Set-TextValue $ws.Range('E41') '  -8.52%  '

Set-TextValue $ws.Range('D42') '497.23'
Set-TextValue $ws.Range('E42') '  -7.86%  '

Set-TextValue $ws.Range('D43') '24.50'
Set-TextValue $ws.Range('E43') '  -1.11%  '

Set-TextValue $ws.Range('D44') '3.64'
Set-TextValue $ws.Range('E44') '  -2.79%  '

Set-TextValue $ws.Range('D45') '0.812'
Set-TextValue $ws.Range('E45') '  -7.00%  '

Set-TextValue $ws.Range('E46') '  -10.36%  '

Set-TextValue $ws.Range('D47') '8.26'
Set-TextValue $ws.Range('E47') '  -5.23%  '

Set-TextValue $ws.Range('B48') 'Filecoin'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D48') '5.32'
Set-TextValue $ws.Range('E48') '  -2.00%  '

Set-TextValue $ws.Range('B49') 'ImmutableX'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D49') '1.60'
Set-TextValue $ws.Range('E49') '  -2.25%  '

Set-TextValue $ws.Range('D50') '52.28'
Set-TextValue $ws.Range('E50') '  -0.42%  '

Set-TextValue $ws.Range('E51') '  -11.46%  '
